# Update cryptos list: refresh prices / 1h volume changes, and fix the
# EthereumClassic / NEARProtocol row ordering (rows 31-32 swapped places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) sometimes hold values that Excel's input parser
# would happily coerce into a Double (e.g. "566.28", "1.00"), which would
# silently drop the original inlineStr/text representation (and its exact
# formatting, like trailing zeros). Force those through as literal text by
# flipping the cell to a text number format before the write, then restore
# the default "Normal" style so no stray formatting is left behind.
function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '59.995.34'
$ws.Range("E2").Value = '  +4.10%  '
$ws.Range("D3").Value = '3.026.00'
$ws.Range("E3").Value = '  +3.10%  '
$ws.Range("E4").Value = '  +0.25%  '
Set-TextValue $ws "D5" '566.28'
$ws.Range("E5").Value = '  +3.03%  '
Set-TextValue $ws "D6" '141.21'
$ws.Range("E6").Value = '  +7.99%  '
$ws.Range("E7").Value = '  -0.07%  '
Set-TextValue $ws "D8" '0.522'
$ws.Range("E8").Value = '  +2.09%  '
$ws.Range("D9").Value = '3.017.28'
$ws.Range("E9").Value = '  +3.02%  '
Set-TextValue $ws "D10" '0.135'
$ws.Range("E10").Value = '  +6.21%  '
$ws.Range("E11").Value = '  +10.98%  '
Set-TextValue $ws "D12" '0.463'
$ws.Range("E12").Value = '  +3.65%  '
$ws.Range("E13").Value = '  +5.41%  '
Set-TextValue $ws "D14" '34.25'
$ws.Range("E14").Value = '  +3.98%  '
$ws.Range("E15").Value = '  +1.95%  '
$ws.Range("D16").Value = '3.527.03'
$ws.Range("E16").Value = '  +3.19%  '
Set-TextValue $ws "D17" '7.23'
$ws.Range("E17").Value = '  +5.40%  '
$ws.Range("D18").Value = '3.023.39'
$ws.Range("E18").Value = '  +3.19%  '
$ws.Range("D19").Value = '59.999.30'
$ws.Range("E19").Value = '  +4.21%  '
Set-TextValue $ws "D20" '439.04'
$ws.Range("E20").Value = '  +4.98%  '
Set-TextValue $ws "D21" '13.75'
$ws.Range("E21").Value = '  +4.24%  '
$ws.Range("E22").Value = '  +5.53%  '
Set-TextValue $ws "D23" '7.16'
$ws.Range("E23").Value = '  +2.44%  '
Set-TextValue $ws "D24" '13.29'
$ws.Range("E24").Value = '  +1.78%  '
Set-TextValue $ws "D25" '81.05'
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  +14.56%  '
Set-TextValue $ws "D28" '1.00'
$ws.Range("E28").Value = '  +0.45%  '
$ws.Range("E29").Value = '  +3.43%  '
Set-TextValue $ws "D30" '7.89'
$ws.Range("E30").Value = '  +5.40%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws "D31" '6.33'
$ws.Range("E31").Value = '  +5.70%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws "D32" '26.16'
$ws.Range("E32").Value = '  +3.60%  '
$ws.Range("E33").Value = '  +4.68%  '
$ws.Range("D34").Value = '0.0₃0792'
$ws.Range("E34").Value = '  +15.83%  '
$ws.Range("E35").Value = '  +7.22%  '
$ws.Range("E36").Value = '  +5.07%  '
$ws.Range("E37").Value = '  +2.89%  '
Set-TextValue $ws "D38" '49.22'
$ws.Range("E38").Value = '  +2.63%  '
Set-TextValue $ws "D39" '8.69'
$ws.Range("E39").Value = '  -0.63%  '
Set-TextValue $ws "D40" '2.80'
$ws.Range("E40").Value = '  +9.54%  '
Set-TextValue $ws "D41" '406.88'
$ws.Range("E41").Value = '  +8.16%  '
$ws.Range("E42").Value = '  +3.06%  '
$ws.Range("D43").Value = '2.789.44'
$ws.Range("E43").Value = '  +3.96%  '
$ws.Range("E44").Value = '  -0.11%  '
Set-TextValue $ws "D45" '0.254'
$ws.Range("E45").Value = '  +6.73%  '
$ws.Range("E46").Value = '  -0.01%  '
Set-TextValue $ws "D47" '123.37'
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("E48").Value = '  +3.41%  '
$ws.Range("E49").Value = '  +1.85%  '
Set-TextValue $ws "D50" '34.20'
$ws.Range("E50").Value = '  +20.46%  '
Set-TextValue $ws "D51" '23.71'
$ws.Range("E51").Value = '  +2.18%  '

